$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-09-01"

# Update the header label cell (I1) text to match the new date
$ws.Range("I1").Value = "2022 (through 09-01)"

# Update September (row 9) running total in column I
$ws.Range("I9").Value = 168

# Add new value for October (row 10) column I
$ws.Range("I10").Value = 3

# Update yearly total (row 14) in column I
$ws.Range("I14").Value = 1142
